$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = 9.053051090061398
$ws.Range("E2").Value = 13.49814074132279
$ws.Range("F2").Value = 34.16091268773713
$ws.Range("G2").Value = 3.623242451708317
$ws.Range("J2").Value = 9.848448925951388
$ws.Range("O2").Value = 25.54137488967321

$ws.Range("D3").Value = 9.062502353047789
$ws.Range("E3").Value = 13.47515106134688
$ws.Range("F3").Value = 33.77373060893166
$ws.Range("G3").Value = 3.626715730940524
$ws.Range("J3").Value = 9.851853114224138
$ws.Range("O3").Value = 25.28327326623468

$ws.Range("D4").Value = 9.069618107795316
$ws.Range("E4").Value = 13.463975922457
$ws.Range("F4").Value = 33.54394766294422
$ws.Range("G4").Value = 3.628960586445144
$ws.Range("J4").Value = 9.855753084105787
$ws.Range("O4").Value = 25.13087928067928

$ws.Range("D5").Value = 9.072847943946005
$ws.Range("E5").Value = 13.46016374882358
$ws.Range("F5").Value = 33.45241244229445
$ws.Range("G5").Value = 3.629903714509452
$ws.Range("J5").Value = 9.857797311550794
$ws.Range("O5").Value = 25.07037136719902

$ws.Range("D6").Value = 9.073404193405823
$ws.Range("E6").Value = 13.45957559886175
$ws.Range("F6").Value = 33.43734307526295
$ws.Range("G6").Value = 3.630062034442503
$ws.Range("J6").Value = 9.858164225332509
$ws.Range("O6").Value = 25.06042215919873

$ws.Range("D7").Value = 9.069660329893461
$ws.Range("E7").Value = 13.46392150404267
$ws.Range("F7").Value = 33.54270453687657
$ws.Range("G7").Value = 3.628973190948789
$ws.Range("J7").Value = 9.855778811462393
$ws.Range("O7").Value = 25.13005671365274

$ws.Range("D8").Value = 9.05603748192066
$ws.Range("E8").Value = 13.48960515884913
$ws.Range("F8").Value = 34.02582672660022
$ws.Range("G8").Value = 3.624416804934584
$ws.Range("J8").Value = 9.849247038929464
$ws.Range("O8").Value = 25.45115907981998

$ws.Range("D9").Value = 9.039735604397226
$ws.Range("E9").Value = 13.56317490237909
$ws.Range("F9").Value = 35.03149178840453
$ws.Range("G9").Value = 3.616367576087087
$ws.Range("J9").Value = 9.850798242082538
$ws.Range("O9").Value = 26.12609398011729

$ws.Range("D10").Value = 9.034100795024408
$ws.Range("E10").Value = 13.63116400232298
$ws.Range("F10").Value = 35.79921196626661
$ws.Range("G10").Value = 3.610987116409229
$ws.Range("J10").Value = 9.860684287964897
$ws.Range("O10").Value = 26.64535503308258

$ws.Range("D11").Value = 9.03291216000437
$ws.Range("E11").Value = 13.66506339699854
$ws.Range("F11").Value = 36.1532782404719
$ws.Range("G11").Value = 3.608653775659333
$ws.Range("J11").Value = 9.867075651012993
$ws.Range("O11").Value = 26.88572683328334

$ws.Range("D12").Value = 9.032659395609095
$ws.Range("E12").Value = 13.67832158929044
$ws.Range("F12").Value = 36.28792533965274
$ws.Range("G12").Value = 3.607786520170719
$ws.Range("J12").Value = 9.869767470888247
$ws.Range("O12").Value = 26.97726771910311

$ws.Range("D13").Value = 9.032705061912399
$ws.Range("E13").Value = 13.67544756740461
$ws.Range("F13").Value = 36.25890330013041
$ws.Range("G13").Value = 3.607972574463286
$ws.Range("J13").Value = 9.869175677268132
$ws.Range("O13").Value = 26.95753102762606

$ws.Range("D14").Value = 9.032887411950668
$ws.Range("E14").Value = 13.66614575069772
$ws.Range("F14").Value = 36.1643450253973
$ws.Range("G14").Value = 3.608582099321673
$ws.Range("J14").Value = 9.867291671810767
$ws.Range("O14").Value = 26.89324805125596

$ws.Range("D15").Value = 9.033024795923772
$ws.Range("E15").Value = 13.66050278847523
$ws.Range("F15").Value = 36.10649589555681
$ws.Range("G15").Value = 3.60895757458505
$ws.Range("J15").Value = 9.866173000034957
$ws.Range("O15").Value = 26.85393784757931

$ws.Range("D16").Value = 9.034206110330613
$ws.Range("E16").Value = 13.62900788030516
$ws.Range("F16").Value = 35.77615939628461
$ws.Range("G16").Value = 3.611141896301158
$ws.Range("J16").Value = 9.860304646131771
$ws.Range("O16").Value = 26.62972294951099

$ws.Range("D17").Value = 9.035282680972772
$ws.Range("E17").Value = 13.6104436073156
$ws.Range("F17").Value = 35.57465633020652
$ws.Range("G17").Value = 3.612511099625995
$ws.Range("J17").Value = 9.85718921014708
$ws.Range("O17").Value = 26.49318197498341

$ws.Range("D18").Value = 9.036031310257808
$ws.Range("E18").Value = 13.60004595055619
$ws.Range("F18").Value = 35.45921881065974
$ws.Range("G18").Value = 3.613309389561933
$ws.Range("J18").Value = 9.855575641806704
$ws.Range("O18").Value = 26.41504338001172

$ws.Range("D19").Value = 9.036307019224624
$ws.Range("E19").Value = 13.59657374906443
$ws.Range("F19").Value = 35.42021674247751
$ws.Range("G19").Value = 3.613581527980637
$ws.Range("J19").Value = 9.855059968827739
$ws.Range("O19").Value = 26.38865748438569

$ws.Range("D20").Value = 9.035154686480707
$ws.Range("E20").Value = 13.61239086607601
$ws.Range("F20").Value = 35.59605979171576
$ws.Range("G20").Value = 3.612364232665131
$ws.Range("J20").Value = 9.857502400698968
$ws.Range("O20").Value = 26.50767657635952

$ws.Range("D21").Value = 9.032828498380541
$ws.Range("E21").Value = 13.66886653821913
$ws.Range("F21").Value = 36.19210460187967
$ws.Range("G21").Value = 3.60840262469864
$ws.Range("J21").Value = 9.86783768764054
$ws.Range("O21").Value = 26.91211612260948

$ws.Range("D22").Value = 9.032458281489083
$ws.Range("E22").Value = 13.70822871588961
$ws.Range("F22").Value = 36.58491920749571
$ws.Range("G22").Value = 3.6059086225994
$ws.Range("J22").Value = 9.876174675601836
$ws.Range("O22").Value = 27.17941864542273

$ws.Range("D23").Value = 9.03255076728813
$ws.Range("E23").Value = 13.68699813834027
$ws.Range("F23").Value = 36.37500870070382
$ws.Range("G23").Value = 3.607231046284066
$ws.Range("J23").Value = 9.871580611595846
$ws.Range("O23").Value = 27.0365083911839

$ws.Range("D24").Value = 9.035212148792517
$ws.Range("E24").Value = 13.61150965277872
$ws.Range("F24").Value = 35.58638200237991
$ws.Range("G24").Value = 3.612430596582852
$ws.Range("J24").Value = 9.85736025418865
$ws.Range("O24").Value = 26.50112243756121

$ws.Range("D25").Value = 9.043031288583045
$ws.Range("E25").Value = 13.54080724258624
$ws.Range("F25").Value = 34.75387368245422
$ws.Range("G25").Value = 3.618450967720438
$ws.Range("J25").Value = 9.848841359131448
$ws.Range("O25").Value = 25.93908695842777
